# Append the 2025-10-03 allocation row (row 32) to the sheet, mirroring the
# existing rows: Date (as literal text, not an Excel date serial), BTC share,
# KAS share.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Write the date as a formula that evaluates to the literal text
# "10/03/2025", then flatten it to a plain value via copy/paste-special.
# This avoids Excel's automatic "looks like a date" parsing that would
# otherwise turn a directly-typed "10/03/2025" into a date serial number.
$ws.Range("A32").Formula = '="10/03/2025"'
$ws.Range("A32").Copy()
$ws.Range("A32").PasteSpecial(-4163)

$ws.Range("B32").Value = 0.1416568838977773
$ws.Range("C32").Value = 0.8583431161022227
